$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September transaction was recorded, which pushes all the
# September (R/S), and eventually August (P/Q), history down by one row.
# Inserting a whole row at row 35 reproduces this cascading shift for every
# column (A, P, Q, R, S, etc.) in one native operation, then we just need to
# populate the newly inserted row's September detail/date cells.
$ws.Rows.Item(35).Insert()

$ws.Range("R35").Value = "bal axisbank"
$ws.Range("S35").Value = "2024-09-09 12:04:31"
